$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after row 5 (new latitude_API / longitude_API rows),
# pushing the old rows 6-10 (depth, biosat, O2_Ar_ratio, ncp, k) down to 8-12.
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()

# Row 2: datetime_utc -- definition now reflects the data product time
$ws.Range("B2").Value = "Data product UTC date and time"

# Row 3: datetime_utc_matlab -- now a PI-provided UTC date/time, same class/format as row 2
$ws.Range("B3").Value = "PI-provided UTC date and time "
$ws.Range("C3").Value = "Date"
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = "YYYY-MM-DD hh:mm:ss"

# Row 4: latitude -> latitude_matlab (PI provided)
$ws.Range("A4").Value = "latitude_matlab"
$ws.Range("B4").Value = "Latitude of sample event provided by PI"

# Row 5: longitude -> longitude_matlab (PI provided)
$ws.Range("A5").Value = "longitude_matlab"
$ws.Range("B5").Value = "Longitude of sample event provided by PI"

# New row 6: latitude_API
$ws.Range("A6").Value = "latitude_API"
$ws.Range("B6").Value = "Latitude of sample event provided by NES-LTER API"
$ws.Range("C6").Value = "numeric"
$ws.Range("D6").Value = "degree"

# New row 7: longitude_API
$ws.Range("A7").Value = "longitude_API"
$ws.Range("B7").Value = "Longitude of sample event provided by NES-LTER API"
$ws.Range("C7").Value = "numeric"
$ws.Range("D7").Value = "degree"

# Column A is now wider to fit the longer attribute names, no longer auto best-fit
$ws.Columns.Item(1).ColumnWidth = 19.3

# Restore the selection Excel leaves behind after the edit
$ws.Range("A4:E7").Select() | Out-Null
